$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New test case row (WAT47) gets appended after the last existing row (37).
# Seed the new row's formatting from row 28, which already carries the
# no-fill / thin-border style combination (border+no wrap on A-D, plain
# border on E) that the new row needs - then overwrite the values.
$ws.Range("A28:E28").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)

$ws.Range("A38").Value = "WAT47"
$ws.Range("C38").Value = 'Verify that when error msg "No Result" is displayed,  Add alternate name should be disabled'
$ws.Range("B38").Value = "WAT-159"
$ws.Range("D38").Value = "Y"

# Move the view so the newly-added row is visible/selected, mirroring the
# author's recorded sheet view after adding the row.
[void]$ws.Activate()
[void]$ws.Range("A1:E38").Select()
